$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Fix up existing values that changed (done first, before row numbers shift) ---
$ws.Range("B12").Value = 385980404667.027
$ws.Range("B16").Value = 7474891657.999
$ws.Range("B87").Value = 494460765005
$ws.Range("B103").Value = 173813766175.705

# --- Step 2: Insert new rows (processed bottom-to-top so row indices stay valid) ---
$ws.Rows.Item(140).Insert()
$ws.Range("A140").Value = 44197
$ws.Range("A140").NumberFormat = "yyyy-mm-dd"
$ws.Range("B140").Value = 3965317044.55
$ws.Range("C140").Value = 'Uganda'
$ws.Range("D140").Value = "World"

$ws.Rows.Item(111).Insert()
$ws.Range("A111").Value = 44197
$ws.Range("A111").NumberFormat = "yyyy-mm-dd"
$ws.Range("B111").Value = 458303194.55
$ws.Range("C111").Value = 'Timor-Leste'
$ws.Range("D111").Value = "World"

$ws.Rows.Item(106).Insert()
$ws.Range("A106").Value = 44197
$ws.Range("A106").NumberFormat = "yyyy-mm-dd"
$ws.Range("B106").Value = 9650541017.488
$ws.Range("C106").Value = 'Papua New Guinea'
$ws.Range("D106").Value = "World"

$ws.Rows.Item(65).Insert()
$ws.Range("A65").Value = 44197
$ws.Range("A65").NumberFormat = "yyyy-mm-dd"
$ws.Range("B65").Value = 15333663013.467
$ws.Range("C65").Value = 'Côte d''Ivoire'
$ws.Range("D65").Value = "World"

$ws.Rows.Item(62).Insert()
$ws.Range("A62").Value = 44197
$ws.Range("A62").NumberFormat = "yyyy-mm-dd"
$ws.Range("B62").Value = 67271404932.41
$ws.Range("C62").Value = 'Iran'
$ws.Range("D62").Value = "World"

$ws.Rows.Item(53).Insert()
$ws.Range("A53").Value = 44197
$ws.Range("A53").NumberFormat = "yyyy-mm-dd"
$ws.Range("B53").Value = 8863588.062
$ws.Range("C53").Value = 'Kiribati'
$ws.Range("D53").Value = "World"

$ws.Rows.Item(49).Insert()
$ws.Range("A49").Value = 44197
$ws.Range("A49").NumberFormat = "yyyy-mm-dd"
$ws.Range("B49").Value = 8031167846.494
$ws.Range("C49").Value = 'Gabon'
$ws.Range("D49").Value = "World"

$ws.Rows.Item(40).Insert()
$ws.Range("A40").Value = 44197
$ws.Range("A40").NumberFormat = "yyyy-mm-dd"
$ws.Range("B40").Value = 14145191.781
$ws.Range("C40").Value = 'Dominica'
$ws.Range("D40").Value = "World"

$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = 44197
$ws.Range("A27").NumberFormat = "yyyy-mm-dd"
$ws.Range("B27").Value = 57725502.927
$ws.Range("C27").Value = 'Central African Rep.'
$ws.Range("D27").Value = "World"

$ws.Rows.Item(26).Insert()
$ws.Range("A26").Value = 44197
$ws.Range("A26").NumberFormat = "yyyy-mm-dd"
$ws.Range("B26").Value = 45366225.058
$ws.Range("C26").Value = 'Cabo Verde'
$ws.Range("D26").Value = "World"

$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = 44197
$ws.Range("A25").NumberFormat = "yyyy-mm-dd"
$ws.Range("B25").Value = 4294055976.357
$ws.Range("C25").Value = 'Cameroon'
$ws.Range("D25").Value = "World"

$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = 44197
$ws.Range("A5").NumberFormat = "yyyy-mm-dd"
$ws.Range("B5").Value = 19147254.49
$ws.Range("C5").Value = 'Antigua and Barbuda'
$ws.Range("D5").Value = "World"
